$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4244
$ws.Range("E2").Value = 548
$ws.Range("F2").Value = 548
$ws.Range("G2").Value = 525
$ws.Range("H2").Value = 386
$ws.Range("I2").Value = 389
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 7947
$ws.Range("L2").Value = 2059
$ws.Range("M2").Value = 5888
$ws.Range("N2").Value = 5861
$ws.Range("O2").Value = 27
$ws.Range("P2").Value = 165
$ws.Range("Q2").Value = 594
$ws.Range("R2").Value = -481
$ws.Range("S2").Value = -32
$ws.Range("T2").Value = 263
$ws.Range("U2").Value = 332
$ws.Range("V2").Value = 249
$ws.Range("W2").Value = 12.92
$ws.Range("X2").Value = 9.1
$ws.Range("Y2").Value = 6.65
$ws.Range("Z2").Value = 4.89
$ws.Range("AA2").Value = 34.97
$ws.Range("AB2").Value = 3323.68
$ws.Range("AC2").Value = 11765
$ws.Range("AD2").Value = 10.15
$ws.Range("AE2").Value = 177208
$ws.Range("AF2").Value = 0.67
$ws.Range("AG2").Value = 1494
$ws.Range("AH2").Value = 1.25
$ws.Range("AI2").Value = $null
$ws.Range("AJ2").Value = 3308710

# Row 3
$ws.Range("D3").Value = 4483
$ws.Range("E3").Value = 580
$ws.Range("F3").Value = 580
$ws.Range("G3").Value = 587
$ws.Range("H3").Value = 443
$ws.Range("I3").Value = 438
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 8093
$ws.Range("L3").Value = 1827
$ws.Range("M3").Value = 6266
$ws.Range("N3").Value = 6223
$ws.Range("O3").Value = 43
$ws.Range("P3").Value = 165
$ws.Range("Q3").Value = 641
$ws.Range("R3").Value = -247
$ws.Range("S3").Value = -198
$ws.Range("T3").Value = 348
$ws.Range("U3").Value = 293
$ws.Range("V3").Value = 97
$ws.Range("W3").Value = 12.94
$ws.Range("X3").Value = 9.890000000000001
$ws.Range("Y3").Value = 7.25
$ws.Range("Z3").Value = 5.53
$ws.Range("AA3").Value = 29.15
$ws.Range("AB3").Value = 3553.94
$ws.Range("AC3").Value = 13236
$ws.Range("AD3").Value = 7.37
$ws.Range("AE3").Value = 188168
$ws.Range("AF3").Value = 0.52
$ws.Range("AG3").Value = 1494
$ws.Range("AH3").Value = 1.53
$ws.Range("AI3").Value = $null
$ws.Range("AJ3").Value = 3308710

# Row 4
$ws.Range("D4").Value = 4557
$ws.Range("E4").Value = 570
$ws.Range("F4").Value = 570
$ws.Range("G4").Value = 542
$ws.Range("H4").Value = 426
$ws.Range("I4").Value = 413
$ws.Range("J4").Value = 13
$ws.Range("K4").Value = 8505
$ws.Range("L4").Value = 1857
$ws.Range("M4").Value = 6648
$ws.Range("N4").Value = 6602
$ws.Range("O4").Value = 46
$ws.Range("P4").Value = 165
$ws.Range("Q4").Value = 724
$ws.Range("R4").Value = -569
$ws.Range("S4").Value = -72
$ws.Range("T4").Value = 445
$ws.Range("U4").Value = 279
$ws.Range("V4").Value = 73
$ws.Range("W4").Value = 12.5
$ws.Range("X4").Value = 9.35
$ws.Range("Y4").Value = 6.44
$ws.Range("Z4").Value = 5.13
$ws.Range("AA4").Value = 27.93
$ws.Range("AB4").Value = 3802.9
$ws.Range("AC4").Value = 12490
$ws.Range("AD4").Value = 5.72
$ws.Range("AE4").Value = 199618
$ws.Range("AF4").Value = 0.36
$ws.Range("AG4").Value = 1494
$ws.Range("AH4").Value = 2.09
$ws.Range("AI4").Value = $null
$ws.Range("AJ4").Value = 3308710

# Row 5
$ws.Range("D5").Value = 4612
$ws.Range("E5").Value = 532
$ws.Range("F5").Value = 532
$ws.Range("G5").Value = 741
$ws.Range("H5").Value = 554
$ws.Range("I5").Value = 552
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 9064
$ws.Range("L5").Value = 2091
$ws.Range("M5").Value = 6973
$ws.Range("N5").Value = 6925
$ws.Range("O5").Value = 48
$ws.Range("P5").Value = 165
$ws.Range("Q5").Value = 568
$ws.Range("R5").Value = -312
$ws.Range("S5").Value = 83
$ws.Range("T5").Value = 381
$ws.Range("U5").Value = 187
$ws.Range("V5").Value = 205
$ws.Range("W5").Value = 11.54
$ws.Range("X5").Value = 12.01
$ws.Range("Y5").Value = 8.16
$ws.Range("Z5").Value = 6.3
$ws.Range("AA5").Value = 29.98
$ws.Range("AB5").Value = 4108.35
$ws.Range("AC5").Value = 16676
$ws.Range("AD5").Value = 6.6
$ws.Range("AE5").Value = 209370
$ws.Range("AF5").Value = 0.53
$ws.Range("AG5").Value = 1494
$ws.Range("AH5").Value = 1.36
$ws.Range("AI5").Value = $null
$ws.Range("AJ5").Value = 3308710

# Row 6
$ws.Range("D6").Value = 8438
$ws.Range("E6").Value = 792
$ws.Range("F6").Value = 792
$ws.Range("G6").Value = 264
$ws.Range("H6").Value = 172
$ws.Range("I6").Value = 153
$ws.Range("K6").Value = 19287
$ws.Range("L6").Value = 11428
$ws.Range("M6").Value = 7859
$ws.Range("N6").Value = 7792
$ws.Range("P6").Value = 195
$ws.Range("Q6").Value = 893
$ws.Range("R6").Value = -3354
$ws.Range("S6").Value = 2587
$ws.Range("T6").Value = 578
$ws.Range("U6").Value = 315
$ws.Range("V6").Value = 7128
$ws.Range("W6").Value = 9.390000000000001
$ws.Range("X6").Value = 2.04
$ws.Range("Y6").Value = 2.08
$ws.Range("Z6").Value = 1.21
$ws.Range("AA6").Value = 145.4
$ws.Range("AB6").Value = 3919.79
$ws.Range("AC6").Value = 4001
$ws.Range("AD6").Value = 25.87
$ws.Range("AE6").Value = 200079
$ws.Range("AF6").Value = 0.52
$ws.Range("AG6").Value = 1500
$ws.Range("AH6").Value = 1.45
$ws.Range("AI6").Value = $null
$ws.Range("AJ6").Value = 3895809

# Row 7
$ws.Range("D7").Value = 8198
$ws.Range("E7").Value = 750
$ws.Range("G7").Value = 386
$ws.Range("H7").Value = 294
$ws.Range("I7").Value = 270
$ws.Range("K7").Value = 19235
$ws.Range("L7").Value = 11050
$ws.Range("M7").Value = 8180
$ws.Range("N7").Value = 8111
$ws.Range("P7").Value = 192
$ws.Range("Q7").Value = 876
$ws.Range("R7").Value = -700
$ws.Range("S7").Value = -432
$ws.Range("T7").Value = 544
$ws.Range("U7").Value = $null
$ws.Range("W7").Value = 9.16
$ws.Range("X7").Value = 3.59
$ws.Range("Y7").Value = 3.39
$ws.Range("Z7").Value = 1.53
$ws.Range("AA7").Value = 135.1
$ws.Range("AC7").Value = 6918
$ws.Range("AD7").Value = 10.34
$ws.Range("AE7").Value = 208269
$ws.Range("AF7").Value = 0.34
$ws.Range("AG7").Value = 1500
$ws.Range("AH7").Value = 2.1
$ws.Range("AI7").Value = 21.68

# Row 8
$ws.Range("D8").Value = 8143
$ws.Range("E8").Value = 805
$ws.Range("G8").Value = 430
$ws.Range("H8").Value = 322
$ws.Range("I8").Value = 322
$ws.Range("K8").Value = 19006
$ws.Range("L8").Value = 10492
$ws.Range("M8").Value = 8518
$ws.Range("N8").Value = 8445
$ws.Range("P8").Value = 192
$ws.Range("Q8").Value = 734
$ws.Range("R8").Value = -440
$ws.Range("S8").Value = -438
$ws.Range("T8").Value = 532
$ws.Range("U8").Value = $null
$ws.Range("W8").Value = 9.890000000000001
$ws.Range("X8").Value = 3.96
$ws.Range("Y8").Value = 3.9
$ws.Range("Z8").Value = 1.69
$ws.Range("AA8").Value = 123.16
$ws.Range("AC8").Value = 8278
$ws.Range("AD8").Value = 8.09
$ws.Range("AE8").Value = 216845
$ws.Range("AF8").Value = 0.31
$ws.Range("AG8").Value = 1500
$ws.Range("AH8").Value = 2.24
$ws.Range("AI8").Value = 18.12

# Row 9
$ws.Range("D9").Value = 8324
$ws.Range("E9").Value = 845
$ws.Range("G9").Value = 468
$ws.Range("H9").Value = 346
$ws.Range("I9").Value = 346
$ws.Range("K9").Value = 19056
$ws.Range("L9").Value = 10174
$ws.Range("M9").Value = 8882
$ws.Range("N9").Value = 8808
$ws.Range("P9").Value = 192
$ws.Range("Q9").Value = 853
$ws.Range("R9").Value = -584
$ws.Range("S9").Value = -358
$ws.Range("T9").Value = 520
$ws.Range("U9").Value = $null
$ws.Range("W9").Value = 10.15
$ws.Range("X9").Value = 4.16
$ws.Range("Y9").Value = 4.02
$ws.Range("Z9").Value = 1.82
$ws.Range("AA9").Value = 114.54
$ws.Range("AC9").Value = 8894
$ws.Range("AD9").Value = 7.53
$ws.Range("AE9").Value = 226179
$ws.Range("AF9").Value = 0.3
$ws.Range("AG9").Value = 1500
$ws.Range("AH9").Value = 2.24
$ws.Range("AI9").Value = 16.86
